$d = $word.ActiveDocument

# Remove the four unneeded list-item paragraphs that followed the blank
# paragraph under the title ("GIT Repository link.", "Short description of
# web application ...", "Sitemap with all pages.", "Some example
# wireframes of pages ..."). Deleting paragraph #3 repeatedly removes each
# one in turn since the following paragraphs shift up into its place.
$d.Paragraphs(3).Range.Delete()
$d.Paragraphs(3).Range.Delete()
$d.Paragraphs(3).Range.Delete()
$d.Paragraphs(3).Range.Delete()

# Mark the inline picture runs as "do not spell check" (w:noProof) for all
# the screenshot/wireframe images except the first one, which already has
# it set.
for ($i = 2; $i -le $d.InlineShapes.Count; $i++) {
    $d.InlineShapes.Item($i).Range.NoProofing = -1
}
